# Scheduled price-refresh: pushes updated Universalis market-board averages
# and recomputed leve profit figures into each job's profit sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 4203.8174
$ws.Range("I15").Value = 4203.8174
$ws.Range("K15").Value = 12611.4522
$ws.Range("M15").Value = -12442.4522

$ws = $wb.Worksheets.Item("ALC")
# Row 39: Riches' Brew
$ws.Range("H39").Value = 235
$ws.Range("I39").Value = 80
$ws.Range("K39").Value = 240
$ws.Range("M39").Value = 56

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 15963.667
$ws.Range("I132").Value = 2558.0925
$ws.Range("J132").Value = 136613.83
$ws.Range("K132").Value = 7674.2775
$ws.Range("L132").Value = 409841.49
$ws.Range("M132").Value = -5144.2775
$ws.Range("N132").Value = -414901.49

$ws = $wb.Worksheets.Item("ALC")
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 5582.9375
$ws.Range("I137").Value = 1765.9474
$ws.Range("J137").Value = 11161.615
$ws.Range("K137").Value = 5297.8422
$ws.Range("L137").Value = 33484.845
$ws.Range("M137").Value = -2747.8422
$ws.Range("N137").Value = -38584.845

$ws = $wb.Worksheets.Item("ALC")
# Row 138: All-night Crafting
$ws.Range("H138").Value = 1882.6061
$ws.Range("I138").Value = 618.95917
$ws.Range("J138").Value = 3120.98
$ws.Range("K138").Value = 1856.87751
$ws.Range("L138").Value = 9362.940000000001
$ws.Range("M138").Value = 3283.12249
$ws.Range("N138").Value = -19642.94

$ws = $wb.Worksheets.Item("ARM")
# Row 23: A Well-rounded Crew
$ws.Range("H23").Value = 71253
$ws.Range("I23").Value = 75006
$ws.Range("J23").Value = 67500
$ws.Range("K23").Value = 75006
$ws.Range("L23").Value = 67500
$ws.Range("M23").Value = -74747
$ws.Range("N23").Value = -68018

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 9239.563
$ws.Range("I32").Value = 8858.423000000001
$ws.Range("J32").Value = 10930.875
$ws.Range("K32").Value = 8858.423000000001
$ws.Range("L32").Value = 10930.875
$ws.Range("M32").Value = -8571.423000000001
$ws.Range("N32").Value = -11504.875

$ws = $wb.Worksheets.Item("ARM")
# Row 37: Get Shirty
$ws.Range("H37").Value = 47783.332
$ws.Range("J37").Value = 47783.332
$ws.Range("L37").Value = 47783.332
$ws.Range("N37").Value = -48329.332

$ws = $wb.Worksheets.Item("ARM")
# Row 41: Skillet Scandal
$ws.Range("H41").Value = 3428
$ws.Range("I41").Value = 3428
$ws.Range("K41").Value = 3428
$ws.Range("M41").Value = -3014

$ws = $wb.Worksheets.Item("ARM")
# Row 44: Very Slow Array
$ws.Range("H44").Value = 36765.332
$ws.Range("J44").Value = 36765.332
$ws.Range("L44").Value = 36765.332
$ws.Range("N44").Value = -37741.332

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1840.2565
$ws.Range("I45").Value = 1480.7391
$ws.Range("J45").Value = 2357.0625
$ws.Range("K45").Value = 1480.7391
$ws.Range("L45").Value = 2357.0625
$ws.Range("M45").Value = -1103.7391
$ws.Range("N45").Value = -3111.0625

$ws = $wb.Worksheets.Item("ARM")
# Row 55: Employee Retention
$ws.Range("H55").Value = 42000
$ws.Range("J55").Value = 42000
$ws.Range("L55").Value = 42000
$ws.Range("N55").Value = -42630

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 1125.8173
$ws.Range("I61").Value = 1012.15875
$ws.Range("J61").Value = 1364.5
$ws.Range("K61").Value = 1012.15875
$ws.Range("L61").Value = 1364.5
$ws.Range("M61").Value = -800.1587500000001
$ws.Range("N61").Value = -1788.5

$ws = $wb.Worksheets.Item("ARM")
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2370.3447
$ws.Range("I63").Value = 2201.5386
$ws.Range("J63").Value = 3833.3333
$ws.Range("K63").Value = 2201.5386
$ws.Range("L63").Value = 3833.3333
$ws.Range("M63").Value = -1515.5386
$ws.Range("N63").Value = -5205.3333

$ws = $wb.Worksheets.Item("ARM")
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2370.3447
$ws.Range("I66").Value = 2201.5386
$ws.Range("J66").Value = 3833.3333
$ws.Range("K66").Value = 11007.693
$ws.Range("L66").Value = 19166.6665
$ws.Range("M66").Value = -7575.692999999999
$ws.Range("N66").Value = -26030.6665

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1848.2297
$ws.Range("I74").Value = 1850.9642
$ws.Range("J74").Value = 1839.7222
$ws.Range("K74").Value = 1850.9642
$ws.Range("L74").Value = 1839.7222
$ws.Range("M74").Value = -976.9641999999999
$ws.Range("N74").Value = -3587.7222

$ws = $wb.Worksheets.Item("ARM")
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1848.2297
$ws.Range("I77").Value = 1850.9642
$ws.Range("J77").Value = 1839.7222
$ws.Range("K77").Value = 9254.821
$ws.Range("L77").Value = 9198.610999999999
$ws.Range("M77").Value = -4886.821
$ws.Range("N77").Value = -17934.611

$ws = $wb.Worksheets.Item("ARM")
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 7144391.5
$ws.Range("I132").Value = 11364716
$ws.Range("J132").Value = 2304.3845
$ws.Range("K132").Value = 34094148
$ws.Range("L132").Value = 6913.1535
$ws.Range("M132").Value = -34091618
$ws.Range("N132").Value = -11973.1535

$ws = $wb.Worksheets.Item("ARM")
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 1125.8173
$ws.Range("I136").Value = 1012.15875
$ws.Range("J136").Value = 1364.5
$ws.Range("K136").Value = 3036.47625
$ws.Range("L136").Value = 4093.5
$ws.Range("M136").Value = -486.4762500000002
$ws.Range("N136").Value = -9193.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 10205.182
$ws.Range("I82").Value = 2042.8334
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 2042.8334
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = -1659.8334
$ws.Range("N82").Value = -20766

$ws = $wb.Worksheets.Item("BSM")
# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 10205.182
$ws.Range("I85").Value = 2042.8334
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 2042.8334
$ws.Range("L85").Value = 20000
$ws.Range("M85").Value = -716.8334
$ws.Range("N85").Value = -22652

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2758.86
$ws.Range("I31").Value = 1148.1471
$ws.Range("J31").Value = 3588.621
$ws.Range("K31").Value = 1148.1471
$ws.Range("L31").Value = 3588.621
$ws.Range("M31").Value = -853.1470999999999
$ws.Range("N31").Value = -4178.621

$ws = $wb.Worksheets.Item("CRP")
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2758.86
$ws.Range("I34").Value = 1148.1471
$ws.Range("J34").Value = 3588.621
$ws.Range("K34").Value = 1148.1471
$ws.Range("L34").Value = 3588.621
$ws.Range("M34").Value = -946.1470999999999
$ws.Range("N34").Value = -3992.621

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1352.7858
$ws.Range("I58").Value = 1005.96155
$ws.Range("J58").Value = 1916.375
$ws.Range("K58").Value = 1005.96155
$ws.Range("L58").Value = 1916.375
$ws.Range("M58").Value = -802.96155
$ws.Range("N58").Value = -2322.375

$ws = $wb.Worksheets.Item("CRP")
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 172527.72
$ws.Range("I122").Value = 400599.66
$ws.Range("J122").Value = 1473.75
$ws.Range("K122").Value = 1201798.98
$ws.Range("L122").Value = 4421.25
$ws.Range("M122").Value = -1199348.98
$ws.Range("N122").Value = -9321.25

$ws = $wb.Worksheets.Item("CRP")
# Row 136: Turali Quality
$ws.Range("H136").Value = 1352.7858
$ws.Range("I136").Value = 1005.96155
$ws.Range("J136").Value = 1916.375
$ws.Range("K136").Value = 3017.88465
$ws.Range("L136").Value = 5749.125
$ws.Range("M136").Value = -467.88465
$ws.Range("N136").Value = -10849.125

$ws = $wb.Worksheets.Item("CUL")
# Row 31: Food Fight
$ws.Range("H31").Value = 2000
$ws.Range("J31").Value = 2000
$ws.Range("L31").Value = 6000
$ws.Range("N31").Value = -6576

$ws = $wb.Worksheets.Item("CUL")
# Row 137: Creative Chocolate
$ws.Range("H137").Value = 33341362
$ws.Range("J137").Value = 71442420
$ws.Range("L137").Value = 214327260
$ws.Range("N137").Value = -214337460

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 10121.23
$ws.Range("J126").Value = 1831.25
$ws.Range("L126").Value = 5493.75
$ws.Range("N126").Value = -10433.75

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 2579.2856
$ws.Range("I132").Value = 1674.619
$ws.Range("K132").Value = 5023.857
$ws.Range("M132").Value = -2493.857

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 21000.445
$ws.Range("J2").Value = 60001.332
$ws.Range("L2").Value = 60001.332
$ws.Range("N2").Value = -60225.332

$ws = $wb.Worksheets.Item("LTW")
# Row 81: I Need Your Glove Tonight
$ws.Range("H81").Value = 32181
$ws.Range("J81").Value = 32181
$ws.Range("L81").Value = 32181
$ws.Range("N81").Value = -34177

$ws = $wb.Worksheets.Item("LTW")
# Row 84: Halonic Drake Handlers (L)
$ws.Range("H84").Value = 32181
$ws.Range("J84").Value = 32181
$ws.Range("L84").Value = 96543
$ws.Range("N84").Value = -106527

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 1946
$ws.Range("I132").Value = 1207.1282
$ws.Range("K132").Value = 3621.3846
$ws.Range("M132").Value = -1091.3846

$ws = $wb.Worksheets.Item("LTW")
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1214.4667
$ws.Range("I136").Value = 984.0577
$ws.Range("J136").Value = 2712.125
$ws.Range("K136").Value = 2952.1731
$ws.Range("L136").Value = 8136.375
$ws.Range("M136").Value = -402.1731
$ws.Range("N136").Value = -13236.375

$ws = $wb.Worksheets.Item("WVR")
# Row 60: And a Haircut Wouldn't Hurt
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2264866.5
$ws.Range("I126").Value = 5884433
$ws.Range("J126").Value = 2637.375
$ws.Range("K126").Value = 17653299
$ws.Range("L126").Value = 7912.125
$ws.Range("M126").Value = -17650829
$ws.Range("N126").Value = -12852.125

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1578.1571
$ws.Range("I132").Value = 1509.5103
$ws.Range("J132").Value = 1738.3334
$ws.Range("K132").Value = 4528.5309
$ws.Range("L132").Value = 5215.0002
$ws.Range("M132").Value = -1998.5309
$ws.Range("N132").Value = -10275.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 13441.873
$ws.Range("I136").Value = 22790.422
$ws.Range("J136").Value = 1068.7941
$ws.Range("K136").Value = 68371.266
$ws.Range("L136").Value = 3206.3823
$ws.Range("M136").Value = -65821.266
$ws.Range("N136").Value = -8306.382300000001
